$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the gender recode: swap the two DecValTok numeric literals so that
#    `1` = 0, `2` = 1   becomes   `1` = 1, `2` = 0
#    We locate the unique anchor "(gender, " and then narrow in on the two
#    " = N" numeric tokens so that only the digit characters are replaced,
#    preserving each run's distinct style (StringTok / DataTypeTok / NormalTok
#    / DecValTok) instead of collapsing the whole phrase into one run.
# ---------------------------------------------------------------------------
$anchor = $d.Content.Duplicate
$anchorFound = $anchor.Find.Execute("(gender, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $anchorFound) {
    throw "Could not find the gender recode anchor '(gender, '"
}

$scanStart = $anchor.End
$scanEnd = $scanStart + 60

$firstTok = $d.Range($scanStart, $scanEnd)
$firstFound = $firstTok.Find.Execute(" = 0,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $firstFound) {
    throw "Could not find ' = 0,' inside the gender recode block"
}
$firstDigit = $d.Range($firstTok.Start + 3, $firstTok.Start + 4)

$secondTok = $d.Range($scanStart, $scanEnd)
$secondFound = $secondTok.Find.Execute(" = 1)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $secondFound) {
    throw "Could not find ' = 1)' inside the gender recode block"
}
$secondDigit = $d.Range($secondTok.Start + 3, $secondTok.Start + 4)

$firstDigit.Text = "1"
$secondDigit.Text = "0"

# ---------------------------------------------------------------------------
# 2) Flip the "gender" column values in the printed tibble rows.
#    Each row lives in its own single run, so a straightforward
#    Find/Replace on the full line text is safe and keeps styling intact.
# ---------------------------------------------------------------------------
$tableRows = @(
    @{ Old = "##  1 CAN           492.  0.93       NA         NA      0          NA            3"; New = "##  1 CAN           492.  0.93       NA         NA      1          NA            3" },
    @{ Old = "##  2 CAN           394. -0.78        0         NA      0          NA            3"; New = "##  2 CAN           394. -0.78        0         NA      1          NA            3" },
    @{ Old = "##  3 CAN           390. -1.3         0          1      1           2            2"; New = "##  3 CAN           390. -1.3         0          1      0           2            2" },
    @{ Old = "##  4 CAN           504.  0.56        0          2      0           2            3"; New = "##  4 CAN           504.  0.56        0          2      1           2            3" },
    @{ Old = "##  5 CAN           466. -0.03        0          3      1           1           NA"; New = "##  5 CAN           466. -0.03        0          3      0           1           NA" },
    @{ Old = "##  6 CAN           398.  0.74        0          1      0           2            2"; New = "##  6 CAN           398.  0.74        0          1      1           2            2" },
    @{ Old = "##  7 CAN           404. NA          NA         NA      0          NA           NA"; New = "##  7 CAN           404. NA          NA         NA      1          NA           NA" },
    @{ Old = "##  8 CAN           406. -2.58        0          4      0           2           NA"; New = "##  8 CAN           406. -2.58        0          4      1           2           NA" },
    @{ Old = "##  9 CAN           609.  0.88        0          4      1           1           NA"; New = "##  9 CAN           609.  0.88        0          4      0           1           NA" },
    @{ Old = "## 10 CAN           452.  0.44        0          1      0           2           NA"; New = "## 10 CAN           452.  0.44        0          1      1           2           NA" }
)

foreach ($row in $tableRows) {
    $found = $d.Content.Find.Execute($row.Old, $true, $false, $false, $false, $false, $true, 1, $false, $row.New, 2)
    if (-not $found) {
        throw "Could not find table row: $($row.Old)"
    }
}

# ---------------------------------------------------------------------------
# 3) Update the regression coefficient summary table.
# ---------------------------------------------------------------------------
$coefRows = @(
    @{ Old = "## (Intercept) 449.3502     1.1874 378.429  < 2e-16 ***"; New = "## (Intercept) 456.6267     1.2423 367.579  < 2e-16 ***" },
    @{ Old = "## gender        7.2765     0.8016   9.077  < 2e-16 ***"; New = "## gender       -7.2765     0.8016  -9.077  < 2e-16 ***" }
)

foreach ($row in $coefRows) {
    $found = $d.Content.Find.Execute($row.Old, $true, $false, $false, $false, $false, $true, 1, $false, $row.New, 2)
    if (-not $found) {
        throw "Could not find coefficient row: $($row.Old)"
    }
}

# ---------------------------------------------------------------------------
# 4) Update the confidence interval table.
# ---------------------------------------------------------------------------
$ciRows = @(
    @{ Old = "## (Intercept) 447.022873 451.677573"; New = "## (Intercept) 454.191858 459.061551" },
    @{ Old = "## gender        5.705252   8.847711"; New = "## gender       -8.847711  -5.705252" }
)

foreach ($row in $ciRows) {
    $found = $d.Content.Find.Execute($row.Old, $true, $false, $false, $false, $false, $true, 1, $false, $row.New, 2)
    if (-not $found) {
        throw "Could not find confidence interval row: $($row.Old)"
    }
}

Write-Host "All edits applied successfully"
